{"js": "// Update the \"Iterazio amaiera\" (iteration end) date in the first table\n// of the document from 2021/03/14 to 2021/03/15.\nconst body = context.document.body;\n\nconst results = body.search(\"2021/03/14\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace just the match text; this preserves the run's existing\n  // formatting (e.g. the eu-ES language mark) because Word.js applies\n  // the replacement text using the formatting already on the range.\n  results.items[0].insertText(\"2021/03/15\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the \"Iterazio amaiera\" (iteration end) date in the first table\n# of the document from 2021/03/14 to 2021/03/15.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"2021/03/14\",  # FindText\n    $false,        # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap (wdFindContinue)\n    $false,        # Format\n    \"2021/03/15\",  # ReplaceWith\n    2              # Replace (wdReplaceAll)\n)\n"}
